# Disponibilidad.xlsx update — "Actualizar 02-05-2021 01-43-18"
#
# 1) Corrects the timestamp stored in D254:D267 (the previous "batch" of
#    availability checks) from 44232.05062243481 to 44232.05062243056.
# 2) Appends a new batch of 14 availability rows (268-281) — one row per
#    monitored service, cycling through the same Name/URL pairs used by
#    every earlier batch — each stamped with the new check time
#    44232.07169250572, including the matching hyperlink in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) fix the stale timestamp on the previous batch (rows 254-267) ----
for ($r = 254; $r -le 267; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.05062243056
}

# --- 2) append the new batch (rows 268-281) ------------------------------
$newStamp = 44232.07169250572

$rows = @(
    @("Odoo",              "https://www.dataintelligence-group.com/"),
    @("Blackbox",          "https://serviciodashboard.azurewebsites.net/"),
    @("PowerBI",           "https://powerbi.microsoft.com/es-es/"),
    @("Dropbox",           "https://www.dropbox.com/"),
    @("Odoo",              "https://dataintelligence.store/"),
    @("GEE",               "https://app-data-i.users.earthengine.app/"),
    @("UtilidadesOdoo",    "https://odooutil.azurewebsites.net/"),
    @("Filtros Dashboard", "https://filtradordashboard.azurewebsites.net/"),
    @("MapStore",          "https://ide.dataintelligence-group.com/mapstore/#/"),
    @("GeoServer",         "https://ide.dataintelligence-group.com/geoserver/web/?0"),
    @("Tomcat",            "https://ide.dataintelligence-group.com/"),
    @("Shiny",             "https://rpubs.com/dataintelligence/"),
    @("Github",            "https://github.com/Sud-Austral/"),
    @("EZ Exporter",       "https://ezexporter.highviewapps.com/exports/export-profile/")
)

$startRow = 268
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r    = $startRow + $i
    $name = $rows[$i][0]
    $url  = $rows[$i][1]

    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $url
    $ws.Cells.Item($r, 3).Value = "Disponible"
    $ws.Cells.Item($r, 4).Value = $newStamp
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # split off a trailing "#/" fragment the same way the workbook's
    # earlier MapStore hyperlinks do (Address without the fragment,
    # SubAddress carrying the "/" anchor)
    if ($url.Contains("#")) {
        $hashIdx  = $url.IndexOf("#")
        $address  = $url.Substring(0, $hashIdx)
        $location = $url.Substring($hashIdx + 1)
        $ws.Hyperlinks.Add($ws.Range("B$r"), $address, $location)
    } else {
        $ws.Hyperlinks.Add($ws.Range("B$r"), $url)
    }
    $ws.Cells.Item($r, 2).Style = "Hyperlink"
}
